$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19 (shifts existing rows 19+ down by one)
$ws.Rows("19:19").Insert()

# Populate the new row with the "j" -> "Toggle Playback Events" shortcut
$ws.Range("A19").Value = "j"
$ws.Range("B19").Value = "Toggle Playback Events"

# Update selection/view to match target state
$ws.Range("B19").Select()
